$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Los animales saben comunicarse entre si de una manera que hace que otras criaturas se mantengan alertas ante posibles problemas. El miedo se comunica con facilidad en la naturaleza porque para todos es beneficioso ser conscientes de la necesidad de tomar inmediatamente medidas evasivas. Los humanos conservamos la capacidad primitiva de transmitir esta emocion. Lo cual es estupendo... a menos que nos aqueje un falso temor. Por eso cobran fuerza las ansiedades infundadas. Parece que a la gente de tu mundo le preocupa mucho una situacion determinada. Sin embargo, eso no es razon para dejar que te preocupe.'
$ws.Range("B3").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Aunque se nos anima a que hagamos uso de nuestra imaginacion, acabamos tomando prestada la de los demas. Vemos dramas en la television. Vamos al cine. Nos dejamos influir por los anuncios y compramos los productos que han disenado unos desconocidos para nosotros. Lo malo es que cuanto mas dejamos que las visiones de otras personas marquen el tono de nuestra vida, mas limitada acaba siendo nuestra imaginacion. El hecho de que cierta persona no pueda concebir tu idea de lo que es posible no significa que no debas creer en ella. Como tu regente se vincula al innovador Urano, te sientes inspirada. Ten fe en ti.'
$ws.Range("B4").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Es facil pensar que son las diferencias las que causan los conflictos. Pero estos tienden a surgir de las similitudes. Podemos ser poco exigentes con los desconocidos. La gente cuyo aspecto, modo de actuar o de hablar son distintos a los nuestros tiende a fascinarnos, en vez de enfurecernos. Son aquellas personas que proceden de situaciones parecidas a la nuestra - pero que no opinan igual que nosotros, las que nos sacan de quicio. Las batallas mas dramaticas casi siempre se libran por las cosas mas triviales. Hoy puedes tender un puente hacia cierta persona. Las dos teneis algo en comun que hara avanzar vuestros planes.'
$ws.Range("B5").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
?Hay algun modo de que consigas lo que necesitas sin tener que disgustar a nadie? Quieres contentar a todo el mundo pero sabes muy bien lo que quieres. Entonces, ?que deberias hacer? ?En que tienes que ceder para que todo el mundo se sienta integrado en tus objetivos? Necesitas respuestas para poder seguir adelante con cierta sensacion de tranquilidad y control. Pero no necesitas resolverlo todo hoy. Ni tampoco debes sentirte presionada a aceptar las opiniones de los demas sobre lo que hay que hacer. Por algo estas en el asiento del conductor. Toma el control del volante.'
$ws.Range("B6").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Algunas personas no soportan la idea de no saber algo. Por eso, si se les pregunta algo que no saben, se inventan la respuesta. Hacen conjeturas. Y esto seria aceptable si pronunciaran su suposicion con un tinte de duda en la voz. Pero suelen dar su respuesta en un tono autoritario que no admite discusion. Y cuanto mas fanfarronean, mas categoricos se vuelven. A pesar de tu optimismo, alguien ha hecho todo lo posible para convencerte de que cierta situacion no tiene solucion. Estas a punto de descubrir lo contrario.'
$ws.Range("B7").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Cuando nos ponemos en marcha hacia un nuevo destino, es facil que tomemos una ruta equivocada. ?Deberiamos volver sobre nuestros pasos tan pronto como nos hayamos dado cuenta del error? ?O seguir adelante, confiando en que el cosmos tiene buenas razones para desviarnos? Cuando viajamos por carreteras fisicas, Google Maps puede orientarnos. Cuando se trata de viajes interiores, la ayuda es menos obvia y el camino a seguir no es tan facil de adivinar. Sin embargo, tu intuicion es hoy tan fiable como cualquier navegador por satelite. Si confias en ella, veras senales claras de que vas por buen camino.'
$ws.Range("B8").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Estas detectando los planes ocultos de cierta persona. Un escenario no es lo que parece. Esa persona ha disfrazado habilmente sus intenciones y cada vez que piensas o hablas sobre el tema te das cuenta de ello. Y no eres la unica que percibe esta confusion. De hecho, el aire esta tan cargado de suposiciones que podrias cortarlo con un cuchillo. Por eso hay tanta tension y caos. ?Como puedes aliviar esta incomoda sensacion? Al hablar abiertamente sobre lo que se te niega, podras convertir una situacion delicada en una constructiva.'
$ws.Range("B9").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Te estas preparando psicologicamente, para afrontar un reto. ?Por que a veces la vida es tan dificil? ?Para que esperar tensiones cuando deberias poder mirar al futuro con esperanza? Lo que ocurre es que tu sensibilidad natural esta captando las ansiedades de los demas. Se han puesto a la defensiva y su actitud te esta afectando. Pero si piensas en los aspectos positivos y respondes a sus preocupaciones con carino y amabilidad, ?podras influir en su estado de animo? Claro que puedes. Confia en que todo va a salir bien, porque es lo que va a pasar.'
$ws.Range("B10").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Puesto que es mas facil ganar en las discusiones que ganar amigos, a veces es mejor perder en las primeras para conservar los segundos. Lo que no significa que tengas que echarte atras de un asunto que significa mucho para ti. Pero ser consciente de que alguien siente la misma pasion que tu (pero de otra forma) te ayudara a mantener abiertas las lineas de comunicacion. Si estas abierta a lo que venga y eres generosa, podras ver mas alla de los aciertos y errores de esta situacion. Y encontrar un punto de entendimiento que no pondra en peligro tus creencias.'
$ws.Range("B11").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Dicen que debemos tener cuidado con lo que deseamos... y estan en lo cierto. Algunos suenos son mas gratificantes en teoria que en la realidad. Vistos por fuera parecen sofisticados y exoticos. Vistos por dentro, son muy diferentes. Eso no significa que no debamos intentar ir tras ellos. Pero es mejor que lo hagamos sin gafas de color de rosa que oculten la realidad. Necesitas apoyo para poder satisfacer tu ambicion actual. Si tienes las ideas claras y compartes tus planes con una persona clave, podras avanzar mas deprisa de lo que crees.'
$ws.Range("B12").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Si, como dice la cancion, "la vida no es mas que un sueno", ?que ocurre cuando nos despertamos? ?Nos damos la vuelta en la cama y volvemos a dormir, con un vago recuerdo de lo sonado? ?Que es real y que no lo es? Las respuestas a estas preguntas no son tan claras como te gustaria. En estos momentos estas percibiendo como real cierta preocupacion por el futuro, cuando en realidad solo se trata de una fantasia. Conforme se vincula a Venus, en movimiento retrogrado, Urano, tu regente, te proporciona la capacidad de detectar una ilusion desconcertante. Esto te liberara para que tomes esa decision tan dificil de tomar.'
$ws.Range("B13").Value = 'PredicciA3n
                        para hoy                        09 de Agosto de 2023
Atras quedaron los dias en los que se consideraba que aquellas personas que creian en las hadas y la magia eran unas credulas. Avanzando en el tiempo, no hace mucho se consideraba credula a la gente que creia en el calentamiento global. No es mi intencion transmitir una opinion. Simplemente estoy senalando que, sean cuales sean nuestras teorias y creencias, corremos el riesgo de que se burlen de ellas. Y que cuando una persona a la que valoramos critica nuestros puntos de vista, esto puede hacer que dudemos. Hoy no dejes que te convenzan de algo. Tus ideas estan bien pensadas y son validas. Merecen que obres en consecuencia.'
